# [Fonds de solidarite] Add 2020-12-29 data
# Update the "nombre_aides" (C) and "montant_total" (D) figures for the
# rows impacted by the 2020-12-29 data refresh. These columns are stored
# as text in the sheet, so force a Text number format before writing the
# new values to avoid Excel auto-converting them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# Grand Est (reg 44)
Set-TextValue "C36" "787"
Set-TextValue "D36" "3659322.19"

Set-TextValue "C37" "378"
Set-TextValue "D37" "3177413.95"

Set-TextValue "C38" "143"
Set-TextValue "D38" "1829349.71"

# Nouvelle-Aquitaine (reg 75)
Set-TextValue "C92" "470"
Set-TextValue "D92" "4294119.16"

Set-TextValue "C93" "190"
Set-TextValue "D93" "2221473.63"

Set-TextValue "C96" "43"
Set-TextValue "D96" "171991.00"
